$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 678.5143
$ws.Range("I33").Value = 478.84
$ws.Range("J33").Value = 1177.7
$ws.Range("K33").Value = 478.84
$ws.Range("L33").Value = 1177.7
$ws.Range("M33").Value = -249.84
$ws.Range("N33").Value = -1635.7

$ws.Range("H43").Value = 87688.62
$ws.Range("I43").Value = 3816.6667
$ws.Range("J43").Value = 112850.2
$ws.Range("K43").Value = 3816.6667
$ws.Range("L43").Value = 112850.2
$ws.Range("M43").Value = -3747.6667
$ws.Range("N43").Value = -112988.2

$ws.Range("H112").Value = 1783.9131
$ws.Range("J112").Value = 2044.4445
$ws.Range("L112").Value = 6133.333500000001
$ws.Range("N112").Value = -8349.333500000001

$ws.Range("H113").Value = 4648.75
$ws.Range("I113").Value = 3846.6667
$ws.Range("J113").Value = 5130
$ws.Range("K113").Value = 3846.6667
$ws.Range("L113").Value = 5130
$ws.Range("M113").Value = -592.6667000000002
$ws.Range("N113").Value = -11638

$ws.Range("H115").Value = 1400
$ws.Range("I115").Value = 800
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 2400
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -833
$ws.Range("N115").Value = -9134

$ws.Range("H116").Value = 70133.56
$ws.Range("I116").Value = 99441.55
$ws.Range("J116").Value = 5656
$ws.Range("K116").Value = 99441.55
$ws.Range("L116").Value = 5656
$ws.Range("M116").Value = -95999.55
$ws.Range("N116").Value = -12540

$ws.Range("H132").Value = 2810.3635
$ws.Range("I132").Value = 1614.8125
$ws.Range("J132").Value = 8696.154
$ws.Range("K132").Value = 4844.4375
$ws.Range("L132").Value = 26088.462
$ws.Range("M132").Value = -2314.4375
$ws.Range("N132").Value = -31148.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17000
$ws.Range("J24").Value = 17000
$ws.Range("L24").Value = 17000
$ws.Range("N24").Value = -17748

$ws.Range("H32").Value = 18946.053
$ws.Range("I32").Value = 9181.712
$ws.Range("J32").Value = 50951.39
$ws.Range("K32").Value = 9181.712
$ws.Range("L32").Value = 50951.39
$ws.Range("M32").Value = -8894.712
$ws.Range("N32").Value = -51525.39

$ws.Range("H61").Value = 543657.4
$ws.Range("I61").Value = 590591.4399999999
$ws.Range("J61").Value = 503763.5
$ws.Range("K61").Value = 590591.4399999999
$ws.Range("L61").Value = 503763.5
$ws.Range("M61").Value = -590379.4399999999
$ws.Range("N61").Value = -504187.5

$ws.Range("H74").Value = 163569.39
$ws.Range("I74").Value = 205294.25
$ws.Range("J74").Value = 55963.21
$ws.Range("K74").Value = 205294.25
$ws.Range("L74").Value = 55963.21
$ws.Range("M74").Value = -204420.25
$ws.Range("N74").Value = -57711.21

$ws.Range("H77").Value = 163569.39
$ws.Range("I77").Value = 205294.25
$ws.Range("J77").Value = 55963.21
$ws.Range("K77").Value = 1026471.25
$ws.Range("L77").Value = 279816.05
$ws.Range("M77").Value = -1022103.25
$ws.Range("N77").Value = -288552.05

$ws.Range("H100").Value = 17000
$ws.Range("J100").Value = 17000
$ws.Range("L100").Value = 17000
$ws.Range("N100").Value = -19164

$ws.Range("H102").Value = 3495.3125
$ws.Range("I102").Value = 1395
$ws.Range("J102").Value = 18197.5
$ws.Range("K102").Value = 1395
$ws.Range("L102").Value = 18197.5
$ws.Range("M102").Value = 227
$ws.Range("N102").Value = -21441.5

$ws.Range("H136").Value = 543657.4
$ws.Range("I136").Value = 590591.4399999999
$ws.Range("J136").Value = 503763.5
$ws.Range("K136").Value = 1771774.32
$ws.Range("L136").Value = 1511290.5
$ws.Range("M136").Value = -1769224.32
$ws.Range("N136").Value = -1516390.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 610.7692
$ws.Range("J80").Value = 667.8182
$ws.Range("L80").Value = 667.8182
$ws.Range("N80").Value = -2663.8182

$ws.Range("H83").Value = 610.7692
$ws.Range("J83").Value = 667.8182
$ws.Range("L83").Value = 3339.091
$ws.Range("N83").Value = -13323.091

$ws.Range("H107").Value = 1160
$ws.Range("I107").Value = 1380
$ws.Range("J107").Value = 632
$ws.Range("K107").Value = 1380
$ws.Range("L107").Value = 632
$ws.Range("M107").Value = 540
$ws.Range("N107").Value = -4472

$ws.Range("H134").Value = 27944.4
$ws.Range("I134").Value = 38067.11
$ws.Range("J134").Value = 6920.3076
$ws.Range("K134").Value = 114201.33
$ws.Range("L134").Value = 20760.9228
$ws.Range("M134").Value = -111666.33
$ws.Range("N134").Value = -25830.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4312.1934
$ws.Range("I86").Value = 3056.75
$ws.Range("J86").Value = 5651.3335
$ws.Range("K86").Value = 3056.75
$ws.Range("L86").Value = 5651.3335
$ws.Range("M86").Value = -1933.75
$ws.Range("N86").Value = -7897.3335

$ws.Range("H89").Value = 4312.1934
$ws.Range("I89").Value = 3056.75
$ws.Range("J89").Value = 5651.3335
$ws.Range("K89").Value = 15283.75
$ws.Range("L89").Value = 28256.6675
$ws.Range("M89").Value = -9667.75
$ws.Range("N89").Value = -39488.6675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 758.7273
$ws.Range("I5").Value = 522.75
$ws.Range("J5").Value = 980.82355
$ws.Range("K5").Value = 1568.25
$ws.Range("L5").Value = 2942.47065
$ws.Range("M5").Value = -1456.25
$ws.Range("N5").Value = -3166.47065

$ws.Range("H92").Value = 530.4286
$ws.Range("J92").Value = 701
$ws.Range("L92").Value = 2103
$ws.Range("N92").Value = -4599

$ws.Range("H121").Value = 7057834
$ws.Range("I121").Value = 71443580
$ws.Range("J121").Value = 3270437.2
$ws.Range("K121").Value = 214330740
$ws.Range("L121").Value = 9811311.600000001
$ws.Range("M121").Value = -214329430
$ws.Range("N121").Value = -9813931.600000001

$ws.Range("H131").Value = 2640.353
$ws.Range("I131").Value = 3223.5454
$ws.Range("J131").Value = 2361.4348
$ws.Range("K131").Value = 9670.636200000001
$ws.Range("L131").Value = 7084.3044
$ws.Range("M131").Value = -4630.636200000001
$ws.Range("N131").Value = -17164.3044

$ws.Range("H135").Value = 758.7273
$ws.Range("I135").Value = 522.75
$ws.Range("J135").Value = 980.82355
$ws.Range("K135").Value = 4704.75
$ws.Range("L135").Value = 8827.41195
$ws.Range("M135").Value = -2169.75
$ws.Range("N135").Value = -13897.41195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.642857
$ws.Range("I2").Value = 22
$ws.Range("J2").Value = 35.444443
$ws.Range("K2").Value = 22
$ws.Range("L2").Value = 35.444443
$ws.Range("M2").Value = 91
$ws.Range("N2").Value = -261.444443

$ws.Range("H80").Value = 8536.615
$ws.Range("I80").Value = 12142
$ws.Range("J80").Value = 6283.25
$ws.Range("K80").Value = 12142
$ws.Range("L80").Value = 6283.25
$ws.Range("M80").Value = -11144
$ws.Range("N80").Value = -8279.25

$ws.Range("H83").Value = 8536.615
$ws.Range("I83").Value = 12142
$ws.Range("J83").Value = 6283.25
$ws.Range("K83").Value = 60710
$ws.Range("L83").Value = 31416.25
$ws.Range("M83").Value = -55718
$ws.Range("N83").Value = -41400.25

$ws.Range("H132").Value = 3903.0513
$ws.Range("I132").Value = 4160.35
$ws.Range("J132").Value = 3632.2104
$ws.Range("K132").Value = 12481.05
$ws.Range("L132").Value = 10896.6312
$ws.Range("M132").Value = -9951.050000000001
$ws.Range("N132").Value = -15956.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2600.6667
$ws.Range("I61").Value = 2600.6667
$ws.Range("K61").Value = 2600.6667
$ws.Range("M61").Value = -2398.6667

$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622

$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112

$ws.Range("H82").Value = 2848.9167
$ws.Range("I82").Value = 2618.1667
$ws.Range("J82").Value = 3079.6667
$ws.Range("K82").Value = 2618.1667
$ws.Range("L82").Value = 3079.6667
$ws.Range("M82").Value = -2257.1667
$ws.Range("N82").Value = -3801.6667

$ws.Range("H85").Value = 2848.9167
$ws.Range("I85").Value = 2618.1667
$ws.Range("J85").Value = 3079.6667
$ws.Range("K85").Value = 2618.1667
$ws.Range("L85").Value = 3079.6667
$ws.Range("M85").Value = -1370.1667
$ws.Range("N85").Value = -5575.6667

$ws.Range("H113").Value = 2600.6667
$ws.Range("I113").Value = 2600.6667
$ws.Range("K113").Value = 2600.6667
$ws.Range("M113").Value = -430.6667000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1933.7333
$ws.Range("I132").Value = 1067.2778
$ws.Range("J132").Value = 3233.4167
$ws.Range("K132").Value = 3201.8334
$ws.Range("L132").Value = 9700.250100000001
$ws.Range("M132").Value = -671.8334000000004
$ws.Range("N132").Value = -14760.2501

$ws.Range("H136").Value = 235726.81
$ws.Range("I136").Value = 29741.285
$ws.Range("J136").Value = 836517.9399999999
$ws.Range("K136").Value = 89223.855
$ws.Range("L136").Value = 2509553.82
$ws.Range("M136").Value = -86673.855
$ws.Range("N136").Value = -2514653.82
